$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H32").Value = 2417
$ws.Range("I32").Value = 2250
$ws.Range("J32").Value = 2500.5
$ws.Range("K32").Value = 2250
$ws.Range("L32").Value = 2500.5
$ws.Range("M32").Value = -1924
$ws.Range("N32").Value = -3152.5

$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 4017.4443
$ws.Range("I33").Value = 5188.5
$ws.Range("K33").Value = 5188.5
$ws.Range("M33").Value = -4959.5

$ws = $wb.Worksheets.Item(1)
$ws.Range("H43").Value = 4210.643
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 4210.643
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 4210.643
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -4348.643

$ws = $wb.Worksheets.Item(1)
$ws.Range("H69").Value = 11562.272
$ws.Range("I69").Value = 9457.5
$ws.Range("J69").Value = 14088
$ws.Range("K69").Value = 28372.5
$ws.Range("L69").Value = 42264
$ws.Range("M69").Value = -27498.5
$ws.Range("N69").Value = -44012

$ws = $wb.Worksheets.Item(1)
$ws.Range("H70").Value = 3532.3333
$ws.Range("J70").Value = 3738.8
$ws.Range("L70").Value = 11216.4
$ws.Range("N70").Value = -11756.4

$ws = $wb.Worksheets.Item(1)
$ws.Range("H72").Value = 11562.272
$ws.Range("I72").Value = 9457.5
$ws.Range("J72").Value = 14088
$ws.Range("K72").Value = 85117.5
$ws.Range("L72").Value = 126792
$ws.Range("M72").Value = -80749.5
$ws.Range("N72").Value = -135528

$ws = $wb.Worksheets.Item(1)
$ws.Range("H73").Value = 3532.3333
$ws.Range("J73").Value = 3738.8
$ws.Range("L73").Value = 11216.4
$ws.Range("N73").Value = -13088.4

$ws = $wb.Worksheets.Item(1)
$ws.Range("H92").Value = 250644
$ws.Range("I92").Value = 333525.34
$ws.Range("J92").Value = 2000
$ws.Range("K92").Value = 333525.34
$ws.Range("L92").Value = 2000
$ws.Range("M92").Value = -332277.34
$ws.Range("N92").Value = -4496

$ws = $wb.Worksheets.Item(1)
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item(1)
$ws.Range("H98").Value = 2285.4
$ws.Range("I98").Value = 1790.9231
$ws.Range("K98").Value = 1790.9231
$ws.Range("M98").Value = -292.9231

$ws = $wb.Worksheets.Item(1)
$ws.Range("H100").Value = 2439.4211
$ws.Range("I100").Value = 1923.5333
$ws.Range("K100").Value = 1923.5333
$ws.Range("M100").Value = -1382.5333

$ws = $wb.Worksheets.Item(1)
$ws.Range("H111").Value = 677.25
$ws.Range("I111").Value = 631.1429
$ws.Range("K111").Value = 1893.4287
$ws.Range("M111").Value = 1173.5713

$ws = $wb.Worksheets.Item(1)
$ws.Range("H116").Value = 4332.3335
$ws.Range("I116").Value = 4985.8335
$ws.Range("J116").Value = 3025.3333
$ws.Range("K116").Value = 4985.8335
$ws.Range("L116").Value = 3025.3333
$ws.Range("M116").Value = -1543.8335
$ws.Range("N116").Value = -9909.3333

$ws = $wb.Worksheets.Item(1)
$ws.Range("H122").Value = 2285.4
$ws.Range("I122").Value = 1790.9231
$ws.Range("K122").Value = 5372.7693
$ws.Range("M122").Value = -2922.7693

$ws = $wb.Worksheets.Item(1)
$ws.Range("H131").Value = 2308.1428
$ws.Range("J131").Value = 10750
$ws.Range("L131").Value = 32250
$ws.Range("N131").Value = -42330

$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 5491.0586
$ws.Range("I132").Value = 5180.1665
$ws.Range("J132").Value = 6237.2
$ws.Range("K132").Value = 15540.4995
$ws.Range("L132").Value = 18711.6
$ws.Range("M132").Value = -13010.4995
$ws.Range("N132").Value = -23771.6

$ws = $wb.Worksheets.Item(1)
$ws.Range("H137").Value = 35435.438
$ws.Range("I137").Value = 61714.883
$ws.Range("J137").Value = 5652.067
$ws.Range("K137").Value = 185144.649
$ws.Range("L137").Value = 16956.201
$ws.Range("M137").Value = -182594.649
$ws.Range("N137").Value = -22056.201

$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 2123
$ws.Range("I138").Value = 585.913
$ws.Range("J138").Value = 3983.6843
$ws.Range("K138").Value = 1757.739
$ws.Range("L138").Value = 11951.0529
$ws.Range("M138").Value = 3382.261
$ws.Range("N138").Value = -22231.0529

$ws = $wb.Worksheets.Item(1)
$ws.Range("H141").Value = 1382.0869
$ws.Range("I141").Value = 1304.238
$ws.Range("J141").Value = 2199.5
$ws.Range("K141").Value = 3912.714
$ws.Range("L141").Value = 6598.5
$ws.Range("M141").Value = 1267.286
$ws.Range("N141").Value = -16958.5

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 14876.87
$ws.Range("I32").Value = 14808.218
$ws.Range("J32").Value = 15641.857
$ws.Range("K32").Value = 14808.218
$ws.Range("L32").Value = 15641.857
$ws.Range("M32").Value = -14521.218
$ws.Range("N32").Value = -16215.857

$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 10216.538
$ws.Range("I61").Value = 12868.333
$ws.Range("J61").Value = 4250
$ws.Range("K61").Value = 12868.333
$ws.Range("L61").Value = 4250
$ws.Range("M61").Value = -12656.333
$ws.Range("N61").Value = -4674

$ws = $wb.Worksheets.Item(2)
$ws.Range("H63").Value = 2200
$ws.Range("J63").Value = 2200
$ws.Range("L63").Value = 2200
$ws.Range("N63").Value = -3572

$ws = $wb.Worksheets.Item(2)
$ws.Range("H66").Value = 2200
$ws.Range("J66").Value = 2200
$ws.Range("L66").Value = 11000
$ws.Range("N66").Value = -17864

$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 53959.26
$ws.Range("I74").Value = 53959.26
$ws.Range("K74").Value = 53959.26
$ws.Range("M74").Value = -53085.26

$ws = $wb.Worksheets.Item(2)
$ws.Range("H77").Value = 53959.26
$ws.Range("I77").Value = 53959.26
$ws.Range("K77").Value = 269796.3
$ws.Range("M77").Value = -265428.3

$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value = 25362.523
$ws.Range("I132").Value = 28908.764
$ws.Range("K132").Value = 86726.292
$ws.Range("M132").Value = -84196.292

$ws = $wb.Worksheets.Item(2)
$ws.Range("H136").Value = 10216.538
$ws.Range("I136").Value = 12868.333
$ws.Range("J136").Value = 4250
$ws.Range("K136").Value = 38604.999
$ws.Range("L136").Value = 12750
$ws.Range("M136").Value = -36054.999
$ws.Range("N136").Value = -17850

$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 3294
$ws.Range("J20").Value = 4298.8
$ws.Range("L20").Value = 4298.8
$ws.Range("N20").Value = -4792.8

$ws = $wb.Worksheets.Item(3)
$ws.Range("H105").Value = 4048.0417
$ws.Range("I105").Value = 3860.2222
$ws.Range("K105").Value = 3860.2222
$ws.Range("M105").Value = -2113.2222

$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 2488.348
$ws.Range("I107").Value = 1645.4706
$ws.Range("J107").Value = 4876.5
$ws.Range("K107").Value = 1645.4706
$ws.Range("L107").Value = 4876.5
$ws.Range("M107").Value = 274.5293999999999
$ws.Range("N107").Value = -8716.5

$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 3757.1333
$ws.Range("I134").Value = 3654.2593
$ws.Range("K134").Value = 10962.7779
$ws.Range("M134").Value = -8427.777900000001

$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2086
$ws.Range("I31").Value = 1835.0667
$ws.Range("K31").Value = 1835.0667
$ws.Range("M31").Value = -1540.0667

$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 2086
$ws.Range("I34").Value = 1835.0667
$ws.Range("K34").Value = 1835.0667
$ws.Range("M34").Value = -1633.0667

$ws = $wb.Worksheets.Item(4)
$ws.Range("H58").Value = 48015.273
$ws.Range("I58").Value = 64681.375
$ws.Range("K58").Value = 64681.375
$ws.Range("M58").Value = -64478.375

$ws = $wb.Worksheets.Item(4)
$ws.Range("H93").Value = 21687.5
$ws.Range("I93").Value = 21687.5
$ws.Range("K93").Value = 21687.5
$ws.Range("M93").Value = -19815.5

$ws = $wb.Worksheets.Item(4)
$ws.Range("H107").Value = 453.2
$ws.Range("I107").Value = 469.7143
$ws.Range("J107").Value = 222
$ws.Range("K107").Value = 469.7143
$ws.Range("L107").Value = 222
$ws.Range("M107").Value = 1450.2857
$ws.Range("N107").Value = -4062

$ws = $wb.Worksheets.Item(4)
$ws.Range("H122").Value = 2079.3333
$ws.Range("I122").Value = 2079.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6237.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3787.999899999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item(4)
$ws.Range("H132").Value = 2217.3333
$ws.Range("I132").Value = 1951.1538
$ws.Range("K132").Value = 5853.4614
$ws.Range("M132").Value = -3323.4614

$ws = $wb.Worksheets.Item(4)
$ws.Range("H134").Value = 68298.93
$ws.Range("I134").Value = 101547.1
$ws.Range("J134").Value = 1802.6
$ws.Range("K134").Value = 304641.3
$ws.Range("L134").Value = 5407.799999999999
$ws.Range("M134").Value = -302106.3
$ws.Range("N134").Value = -10477.8

$ws = $wb.Worksheets.Item(4)
$ws.Range("H136").Value = 48015.273
$ws.Range("I136").Value = 64681.375
$ws.Range("K136").Value = 194044.125
$ws.Range("M136").Value = -191494.125

$ws = $wb.Worksheets.Item(5)
$ws.Range("H50").Value = 100177.7
$ws.Range("J50").Value = 143069.14
$ws.Range("L50").Value = 429207.42
$ws.Range("N50").Value = -430169.42

$ws = $wb.Worksheets.Item(5)
$ws.Range("H53").Value = 100177.7
$ws.Range("J53").Value = 143069.14
$ws.Range("L53").Value = 429207.42
$ws.Range("N53").Value = -430169.42

$ws = $wb.Worksheets.Item(5)
$ws.Range("H57").Value = 13483.167
$ws.Range("I57").Value = 9966.667
$ws.Range("J57").Value = 16999.666
$ws.Range("K57").Value = 29900.001
$ws.Range("L57").Value = 50998.99800000001
$ws.Range("M57").Value = -29341.001
$ws.Range("N57").Value = -52116.99800000001

$ws = $wb.Worksheets.Item(5)
$ws.Range("H105").Value = 9407
$ws.Range("J105").Value = 9407
$ws.Range("L105").Value = 28221
$ws.Range("N105").Value = -33463

$ws = $wb.Worksheets.Item(6)
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item(6)
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws = $wb.Worksheets.Item(6)
$ws.Range("H12").Value = 30000
$ws.Range("J12").Value = 30000
$ws.Range("L12").Value = 30000
$ws.Range("N12").Value = -30280

$ws = $wb.Worksheets.Item(6)
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30490

$ws = $wb.Worksheets.Item(6)
$ws.Range("H24").Value = 2719.7
$ws.Range("J24").Value = 2719.7
$ws.Range("L24").Value = 2719.7
$ws.Range("N24").Value = -3065.7

$ws = $wb.Worksheets.Item(6)
$ws.Range("H29").Value = 15425
$ws.Range("J29").Value = 15425
$ws.Range("L29").Value = 15425
$ws.Range("N29").Value = -16005

$ws = $wb.Worksheets.Item(6)
$ws.Range("H40").Value = 10000000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 7225
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730

$ws = $wb.Worksheets.Item(6)
$ws.Range("H73").Value = 7225
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064

$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 5516
$ws.Range("I80").Value = 6074.3
$ws.Range("J80").Value = 4399.4
$ws.Range("K80").Value = 6074.3
$ws.Range("L80").Value = 4399.4
$ws.Range("M80").Value = -5076.3
$ws.Range("N80").Value = -6395.4

$ws = $wb.Worksheets.Item(6)
$ws.Range("H83").Value = 5516
$ws.Range("I83").Value = 6074.3
$ws.Range("J83").Value = 4399.4
$ws.Range("K83").Value = 30371.5
$ws.Range("L83").Value = 21997
$ws.Range("M83").Value = -25379.5
$ws.Range("N83").Value = -31981

$ws = $wb.Worksheets.Item(6)
$ws.Range("H126").Value = 6488.9375
$ws.Range("I126").Value = 5895
$ws.Range("J126").Value = 7795.6
$ws.Range("K126").Value = 17685
$ws.Range("L126").Value = 23386.8
$ws.Range("M126").Value = -15215
$ws.Range("N126").Value = -28326.8

$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 32583.158
$ws.Range("I132").Value = 42573.855
$ws.Range("K132").Value = 127721.565
$ws.Range("M132").Value = -125191.565

$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 14073
$ws.Range("I7").Value = 17977.889
$ws.Range("J7").Value = 5287
$ws.Range("K7").Value = 17977.889
$ws.Range("L7").Value = 5287
$ws.Range("M7").Value = -17865.889
$ws.Range("N7").Value = -5511

$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 2499.2778
$ws.Range("I16").Value = 3305.32
$ws.Range("K16").Value = 3305.32
$ws.Range("M16").Value = -3135.32

$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 63987.055
$ws.Range("J22").Value = 2651.4666
$ws.Range("L22").Value = 2651.4666
$ws.Range("N22").Value = -3241.4666

$ws = $wb.Worksheets.Item(7)
$ws.Range("H23").Value = 18000
$ws.Range("I23").Value = 18000
$ws.Range("K23").Value = 18000
$ws.Range("M23").Value = -17770

$ws = $wb.Worksheets.Item(7)
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws = $wb.Worksheets.Item(7)
$ws.Range("H27").Value = 63987.055
$ws.Range("J27").Value = 2651.4666
$ws.Range("L27").Value = 2651.4666
$ws.Range("N27").Value = -2865.4666

$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 13267.529
$ws.Range("I40").Value = 13603.692
$ws.Range("J40").Value = 12175
$ws.Range("K40").Value = 13603.692
$ws.Range("L40").Value = 12175
$ws.Range("M40").Value = -13467.692
$ws.Range("N40").Value = -12447

$ws = $wb.Worksheets.Item(7)
$ws.Range("H42").Value = 39899.5
$ws.Range("J42").Value = 39899.5
$ws.Range("L42").Value = 39899.5
$ws.Range("N42").Value = -41025.5

$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 19055.059
$ws.Range("I46").Value = 32636.875
$ws.Range("K46").Value = 32636.875
$ws.Range("M46").Value = -32448.875

$ws = $wb.Worksheets.Item(7)
$ws.Range("H49").Value = 39899.5
$ws.Range("J49").Value = 39899.5
$ws.Range("L49").Value = 39899.5
$ws.Range("N49").Value = -40193.5

$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 5273.4287
$ws.Range("I61").Value = 5475.5
$ws.Range("J61").Value = 5004
$ws.Range("K61").Value = 5475.5
$ws.Range("L61").Value = 5004
$ws.Range("M61").Value = -5273.5
$ws.Range("N61").Value = -5408

$ws = $wb.Worksheets.Item(7)
$ws.Range("H68").Value = 2703
$ws.Range("I68").Value = 2349.1667
$ws.Range("J68").Value = 2968.375
$ws.Range("K68").Value = 2349.1667
$ws.Range("L68").Value = 2968.375
$ws.Range("M68").Value = -1600.1667
$ws.Range("N68").Value = -4466.375

$ws = $wb.Worksheets.Item(7)
$ws.Range("H71").Value = 2703
$ws.Range("I71").Value = 2349.1667
$ws.Range("J71").Value = 2968.375
$ws.Range("K71").Value = 11745.8335
$ws.Range("L71").Value = 14841.875
$ws.Range("M71").Value = -8001.833500000001
$ws.Range("N71").Value = -22329.875

$ws = $wb.Worksheets.Item(7)
$ws.Range("H100").Value = 4294.5713
$ws.Range("I100").Value = 3993
$ws.Range("K100").Value = 3993
$ws.Range("M100").Value = -3452

$ws = $wb.Worksheets.Item(7)
$ws.Range("H113").Value = 5273.4287
$ws.Range("I113").Value = 5475.5
$ws.Range("J113").Value = 5004
$ws.Range("K113").Value = 5475.5
$ws.Range("L113").Value = 5004
$ws.Range("M113").Value = -3305.5
$ws.Range("N113").Value = -9344

$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 2842.111
$ws.Range("I122").Value = 2226.2856
$ws.Range("J122").Value = 4997.5
$ws.Range("K122").Value = 6678.8568
$ws.Range("L122").Value = 14992.5
$ws.Range("M122").Value = -4228.8568
$ws.Range("N122").Value = -19892.5

$ws = $wb.Worksheets.Item(7)
$ws.Range("H126").Value = 14073
$ws.Range("I126").Value = 17977.889
$ws.Range("J126").Value = 5287
$ws.Range("K126").Value = 53933.667
$ws.Range("L126").Value = 15861
$ws.Range("M126").Value = -51463.667
$ws.Range("N126").Value = -20801

$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 5399.263
$ws.Range("J136").Value = 5999.1665
$ws.Range("L136").Value = 17997.4995
$ws.Range("N136").Value = -23097.4995

$ws = $wb.Worksheets.Item(8)
$ws.Range("H6").Value = 23874.5
$ws.Range("J6").Value = 23874.5
$ws.Range("L6").Value = 23874.5
$ws.Range("N6").Value = -24104.5

$ws = $wb.Worksheets.Item(8)
$ws.Range("H14").Value = 924.375
$ws.Range("I14").Value = 913.7143
$ws.Range("K14").Value = 913.7143
$ws.Range("M14").Value = -745.7143

$ws = $wb.Worksheets.Item(8)
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item(8)
$ws.Range("H20").Value = 49949.5
$ws.Range("J20").Value = 49949.5
$ws.Range("L20").Value = 49949.5
$ws.Range("N20").Value = -50429.5

$ws = $wb.Worksheets.Item(8)
$ws.Range("H26").Value = 11255.5
$ws.Range("J26").Value = 15601.6
$ws.Range("L26").Value = 15601.6
$ws.Range("N26").Value = -16187.6

$ws = $wb.Worksheets.Item(8)
$ws.Range("H30").Value = 5000
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 101833.37
$ws.Range("I62").Value = 7096.6665
$ws.Range("J62").Value = 137359.62
$ws.Range("K62").Value = 7096.6665
$ws.Range("L62").Value = 137359.62
$ws.Range("M62").Value = -6472.6665
$ws.Range("N62").Value = -138607.62

$ws = $wb.Worksheets.Item(8)
$ws.Range("H65").Value = 101833.37
$ws.Range("I65").Value = 7096.6665
$ws.Range("J65").Value = 137359.62
$ws.Range("K65").Value = 35483.3325
$ws.Range("L65").Value = 686798.1
$ws.Range("M65").Value = -32363.3325
$ws.Range("N65").Value = -693038.1

$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 1823
$ws.Range("I81").Value = 1790.5
$ws.Range("K81").Value = 3581
$ws.Range("M81").Value = -2520

$ws = $wb.Worksheets.Item(8)
$ws.Range("H84").Value = 1823
$ws.Range("I84").Value = 1790.5
$ws.Range("K84").Value = 17905
$ws.Range("M84").Value = -12601

$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 25179.479
$ws.Range("I132").Value = 30048.795
$ws.Range("J132").Value = 4079.111
$ws.Range("K132").Value = 90146.385
$ws.Range("L132").Value = 12237.333
$ws.Range("M132").Value = -87616.385
$ws.Range("N132").Value = -17297.333
